$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 16 data: CatBoost benchmark results (added below the existing
# GBDT row, after the LightGBM/linear/FM/GBDT baselines already on the sheet)
$ws.Range("A16").Value = "CatBoost"
$ws.Range("B16").Value = 0.6945
$ws.Range("C16").Value = 0.661
$ws.Range("D16").Value = 0.6544
$ws.Range("E16").Formula = "=1000*(D16-0.6527)"
$ws.Range("H16").Formula = "=D16-C16"

# Update selection to match final state
$ws.Range("H13").Select() | Out-Null
